$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting template from last existing data row (149) for column A and E
$ws.Range("A149").Copy() | Out-Null
$ws.Range("A150:A161").PasteSpecial(-4122) | Out-Null
$ws.Range("E149").Copy() | Out-Null
$ws.Range("E150:E161").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 150 (Indice 149)
$ws.Cells.Item(150, 1).Value = 149
$ws.Cells.Item(150, 2).Value = "belgium"
$ws.Cells.Item(150, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(150, 4).Value = "2023-2024"
$ws.Cells.Item(150, 5).Value = 45283.66666666666
$ws.Cells.Item(150, 6).Value = "Eupen"
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = "Royale Union SG"
$ws.Cells.Item(150, 9).Value = 2
$ws.Cells.Item(150, 10).Value = 8.76
$ws.Cells.Item(150, 11).Value = "17/12/2023 16:12"
$ws.Cells.Item(150, 12).Value = 14.47
$ws.Cells.Item(150, 13).Value = "23/12/2023 15:49"
$ws.Cells.Item(150, 14).Value = 5.72
$ws.Cells.Item(150, 15).Value = "17/12/2023 16:12"
$ws.Cells.Item(150, 16).Value = 7.21
$ws.Cells.Item(150, 17).Value = "23/12/2023 15:49"
$ws.Cells.Item(150, 18).Value = 1.28
$ws.Cells.Item(150, 19).Value = "17/12/2023 16:12"
$ws.Cells.Item(150, 20).Value = 1.2
$ws.Cells.Item(150, 21).Value = "23/12/2023 15:49"
$ws.Cells.Item(150, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/eupen-royale-union-sg/GYgXXpfM/"

# Row 151 (Indice 150)
$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = "belgium"
$ws.Cells.Item(151, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(151, 4).Value = "2023-2024"
$ws.Cells.Item(151, 5).Value = 45283.76041666666
$ws.Cells.Item(151, 6).Value = "Antwerp"
$ws.Cells.Item(151, 7).Value = 2
$ws.Cells.Item(151, 8).Value = "Westerlo"
$ws.Cells.Item(151, 9).Value = 2
$ws.Cells.Item(151, 10).Value = 1.35
$ws.Cells.Item(151, 11).Value = "17/12/2023 18:43"
$ws.Cells.Item(151, 12).Value = 1.32
$ws.Cells.Item(151, 13).Value = "23/12/2023 18:06"
$ws.Cells.Item(151, 14).Value = 5.1
$ws.Cells.Item(151, 15).Value = "17/12/2023 18:43"
$ws.Cells.Item(151, 16).Value = 5.76
$ws.Cells.Item(151, 17).Value = "23/12/2023 18:14"
$ws.Cells.Item(151, 18).Value = 7.33
$ws.Cells.Item(151, 19).Value = "17/12/2023 18:43"
$ws.Cells.Item(151, 20).Value = 9.27
$ws.Cells.Item(151, 21).Value = "23/12/2023 18:14"
$ws.Cells.Item(151, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/antwerp-westerlo/UHXVTGyi/"

# Row 152 (Indice 151)
$ws.Cells.Item(152, 1).Value = 151
$ws.Cells.Item(152, 2).Value = "belgium"
$ws.Cells.Item(152, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(152, 4).Value = "2023-2024"
$ws.Cells.Item(152, 5).Value = 45283.76041666666
$ws.Cells.Item(152, 6).Value = "St. Truiden"
$ws.Cells.Item(152, 7).Value = 1
$ws.Cells.Item(152, 8).Value = "Charleroi"
$ws.Cells.Item(152, 9).Value = 0
$ws.Cells.Item(152, 10).Value = 2.39
$ws.Cells.Item(152, 11).Value = "16/12/2023 21:12"
$ws.Cells.Item(152, 12).Value = 2.4
$ws.Cells.Item(152, 13).Value = "23/12/2023 18:06"
$ws.Cells.Item(152, 14).Value = 3.27
$ws.Cells.Item(152, 15).Value = "16/12/2023 21:12"
$ws.Cells.Item(152, 16).Value = 3.31
$ws.Cells.Item(152, 17).Value = "23/12/2023 18:06"
$ws.Cells.Item(152, 18).Value = 2.96
$ws.Cells.Item(152, 19).Value = "16/12/2023 21:12"
$ws.Cells.Item(152, 20).Value = 3.17
$ws.Cells.Item(152, 21).Value = "23/12/2023 18:06"
$ws.Cells.Item(152, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-truiden-charleroi/ng2GyTfc/"

# Row 153 (Indice 152)
$ws.Cells.Item(153, 1).Value = 152
$ws.Cells.Item(153, 2).Value = "belgium"
$ws.Cells.Item(153, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(153, 4).Value = "2023-2024"
$ws.Cells.Item(153, 5).Value = 45283.86458333334
$ws.Cells.Item(153, 6).Value = "Anderlecht"
$ws.Cells.Item(153, 7).Value = 2
$ws.Cells.Item(153, 8).Value = "Genk"
$ws.Cells.Item(153, 9).Value = 1
$ws.Cells.Item(153, 10).Value = 2.46
$ws.Cells.Item(153, 11).Value = "17/12/2023 19:43"
$ws.Cells.Item(153, 12).Value = 2.75
$ws.Cells.Item(153, 13).Value = "23/12/2023 20:44"
$ws.Cells.Item(153, 14).Value = 3.5
$ws.Cells.Item(153, 15).Value = "17/12/2023 19:43"
$ws.Cells.Item(153, 16).Value = 3.55
$ws.Cells.Item(153, 17).Value = "23/12/2023 20:44"
$ws.Cells.Item(153, 18).Value = 2.71
$ws.Cells.Item(153, 19).Value = "17/12/2023 19:43"
$ws.Cells.Item(153, 20).Value = 2.58
$ws.Cells.Item(153, 21).Value = "23/12/2023 20:44"
$ws.Cells.Item(153, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/anderlecht-genk/tG1KzmA3/"

# Row 154 (Indice 153)
$ws.Cells.Item(154, 1).Value = 153
$ws.Cells.Item(154, 2).Value = "belgium"
$ws.Cells.Item(154, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(154, 4).Value = "2023-2024"
$ws.Cells.Item(154, 5).Value = 45286.5625
$ws.Cells.Item(154, 6).Value = "Genk"
$ws.Cells.Item(154, 7).Value = 3
$ws.Cells.Item(154, 8).Value = "Antwerp"
$ws.Cells.Item(154, 9).Value = 0
$ws.Cells.Item(154, 10).Value = 2.24
$ws.Cells.Item(154, 11).Value = "23/12/2023 21:12"
$ws.Cells.Item(154, 12).Value = 2.37
$ws.Cells.Item(154, 13).Value = "26/12/2023 13:29"
$ws.Cells.Item(154, 14).Value = 3.57
$ws.Cells.Item(154, 15).Value = "23/12/2023 21:12"
$ws.Cells.Item(154, 16).Value = 3.49
$ws.Cells.Item(154, 17).Value = "26/12/2023 13:29"
$ws.Cells.Item(154, 18).Value = 2.97
$ws.Cells.Item(154, 19).Value = "23/12/2023 21:12"
$ws.Cells.Item(154, 20).Value = 3.08
$ws.Cells.Item(154, 21).Value = "26/12/2023 13:29"
$ws.Cells.Item(154, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-antwerp/8WU7RrXk/"

# Row 155 (Indice 154)
$ws.Cells.Item(155, 1).Value = 154
$ws.Cells.Item(155, 2).Value = "belgium"
$ws.Cells.Item(155, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(155, 4).Value = "2023-2024"
$ws.Cells.Item(155, 5).Value = 45286.66666666666
$ws.Cells.Item(155, 6).Value = "Kortrijk"
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = "Gent"
$ws.Cells.Item(155, 9).Value = 2
$ws.Cells.Item(155, 10).Value = 7.2
$ws.Cells.Item(155, 11).Value = "23/12/2023 15:12"
$ws.Cells.Item(155, 12).Value = 10.78
$ws.Cells.Item(155, 13).Value = "26/12/2023 15:48"
$ws.Cells.Item(155, 14).Value = 5.33
$ws.Cells.Item(155, 15).Value = "23/12/2023 15:12"
$ws.Cells.Item(155, 16).Value = 5.93
$ws.Cells.Item(155, 17).Value = "26/12/2023 15:48"
$ws.Cells.Item(155, 18).Value = 1.33
$ws.Cells.Item(155, 19).Value = "23/12/2023 15:12"
$ws.Cells.Item(155, 20).Value = 1.28
$ws.Cells.Item(155, 21).Value = "26/12/2023 15:47"
$ws.Cells.Item(155, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/kortrijk-gent/Q7u3SOHq/"

# Row 156 (Indice 155)
$ws.Cells.Item(156, 1).Value = 155
$ws.Cells.Item(156, 2).Value = "belgium"
$ws.Cells.Item(156, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(156, 4).Value = "2023-2024"
$ws.Cells.Item(156, 5).Value = 45286.66666666666
$ws.Cells.Item(156, 6).Value = "Leuven"
$ws.Cells.Item(156, 7).Value = 3
$ws.Cells.Item(156, 8).Value = "Eupen"
$ws.Cells.Item(156, 9).Value = 0
$ws.Cells.Item(156, 10).Value = 1.75
$ws.Cells.Item(156, 11).Value = "23/12/2023 18:12"
$ws.Cells.Item(156, 12).Value = 1.87
$ws.Cells.Item(156, 13).Value = "26/12/2023 15:59"
$ws.Cells.Item(156, 14).Value = 4.04
$ws.Cells.Item(156, 15).Value = "23/12/2023 18:12"
$ws.Cells.Item(156, 16).Value = 3.99
$ws.Cells.Item(156, 17).Value = "26/12/2023 15:59"
$ws.Cells.Item(156, 18).Value = 3.98
$ws.Cells.Item(156, 19).Value = "23/12/2023 18:12"
$ws.Cells.Item(156, 20).Value = 3.99
$ws.Cells.Item(156, 21).Value = "26/12/2023 15:59"
$ws.Cells.Item(156, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/leuven-eupen/W4G2mQg9/"

# Row 157 (Indice 156)
$ws.Cells.Item(157, 1).Value = 156
$ws.Cells.Item(157, 2).Value = "belgium"
$ws.Cells.Item(157, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(157, 4).Value = "2023-2024"
$ws.Cells.Item(157, 5).Value = 45286.77083333334
$ws.Cells.Item(157, 6).Value = "Westerlo"
$ws.Cells.Item(157, 7).Value = 3
$ws.Cells.Item(157, 8).Value = "RWDM"
$ws.Cells.Item(157, 9).Value = 0
$ws.Cells.Item(157, 10).Value = 1.97
$ws.Cells.Item(157, 11).Value = "23/12/2023 18:42"
$ws.Cells.Item(157, 12).Value = 1.85
$ws.Cells.Item(157, 13).Value = "26/12/2023 18:28"
$ws.Cells.Item(157, 14).Value = 3.64
$ws.Cells.Item(157, 15).Value = "23/12/2023 18:42"
$ws.Cells.Item(157, 16).Value = 3.85
$ws.Cells.Item(157, 17).Value = "26/12/2023 18:28"
$ws.Cells.Item(157, 18).Value = 3.5
$ws.Cells.Item(157, 19).Value = "23/12/2023 18:42"
$ws.Cells.Item(157, 20).Value = 4.26
$ws.Cells.Item(157, 21).Value = "26/12/2023 18:28"
$ws.Cells.Item(157, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/westerlo-rwd-molenbeek/zwNgknPd/"

# Row 158 (Indice 157)
$ws.Cells.Item(158, 1).Value = 157
$ws.Cells.Item(158, 2).Value = "belgium"
$ws.Cells.Item(158, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(158, 4).Value = "2023-2024"
$ws.Cells.Item(158, 5).Value = 45286.86458333334
$ws.Cells.Item(158, 6).Value = "Club Brugge KV"
$ws.Cells.Item(158, 7).Value = 1
$ws.Cells.Item(158, 8).Value = "Royale Union SG"
$ws.Cells.Item(158, 9).Value = 1
$ws.Cells.Item(158, 10).Value = 2.1
$ws.Cells.Item(158, 11).Value = "23/12/2023 18:12"
$ws.Cells.Item(158, 12).Value = 2.3
$ws.Cells.Item(158, 13).Value = "26/12/2023 20:40"
$ws.Cells.Item(158, 14).Value = 3.48
$ws.Cells.Item(158, 15).Value = "23/12/2023 18:12"
$ws.Cells.Item(158, 16).Value = 3.45
$ws.Cells.Item(158, 17).Value = "26/12/2023 20:40"
$ws.Cells.Item(158, 18).Value = 3.34
$ws.Cells.Item(158, 19).Value = "23/12/2023 18:12"
$ws.Cells.Item(158, 20).Value = 3.24
$ws.Cells.Item(158, 21).Value = "26/12/2023 20:40"
$ws.Cells.Item(158, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/club-brugge-royale-union-sg/ETLcl6v3/"

# Row 159 (Indice 158)
$ws.Cells.Item(159, 1).Value = 158
$ws.Cells.Item(159, 2).Value = "belgium"
$ws.Cells.Item(159, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(159, 4).Value = "2023-2024"
$ws.Cells.Item(159, 5).Value = 45287.77083333334
$ws.Cells.Item(159, 6).Value = "St. Liege"
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = "St. Truiden"
$ws.Cells.Item(159, 9).Value = 1
$ws.Cells.Item(159, 10).Value = 2.13
$ws.Cells.Item(159, 11).Value = "23/12/2023 18:42"
$ws.Cells.Item(159, 12).Value = 2.59
$ws.Cells.Item(159, 13).Value = "27/12/2023 18:29"
$ws.Cells.Item(159, 14).Value = 3.34
$ws.Cells.Item(159, 15).Value = "23/12/2023 18:42"
$ws.Cells.Item(159, 16).Value = 3.11
$ws.Cells.Item(159, 17).Value = "27/12/2023 18:29"
$ws.Cells.Item(159, 18).Value = 3.38
$ws.Cells.Item(159, 19).Value = "23/12/2023 18:42"
$ws.Cells.Item(159, 20).Value = 3.08
$ws.Cells.Item(159, 21).Value = "27/12/2023 18:29"
$ws.Cells.Item(159, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-liege-st-truiden/4IEAo4OL/"

# Row 160 (Indice 159)
$ws.Cells.Item(160, 1).Value = 159
$ws.Cells.Item(160, 2).Value = "belgium"
$ws.Cells.Item(160, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(160, 4).Value = "2023-2024"
$ws.Cells.Item(160, 5).Value = 45287.86458333334
$ws.Cells.Item(160, 6).Value = "Charleroi"
$ws.Cells.Item(160, 7).Value = 3
$ws.Cells.Item(160, 8).Value = "KV Mechelen"
$ws.Cells.Item(160, 9).Value = 1
$ws.Cells.Item(160, 10).Value = 2.07
$ws.Cells.Item(160, 11).Value = "23/12/2023 18:42"
$ws.Cells.Item(160, 12).Value = 2.42
$ws.Cells.Item(160, 13).Value = "27/12/2023 20:42"
$ws.Cells.Item(160, 14).Value = 3.44
$ws.Cells.Item(160, 15).Value = "23/12/2023 18:42"
$ws.Cells.Item(160, 16).Value = 3.2
$ws.Cells.Item(160, 17).Value = "27/12/2023 20:42"
$ws.Cells.Item(160, 18).Value = 3.42
$ws.Cells.Item(160, 19).Value = "23/12/2023 18:42"
$ws.Cells.Item(160, 20).Value = 3.26
$ws.Cells.Item(160, 21).Value = "27/12/2023 20:42"
$ws.Cells.Item(160, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/charleroi-kv-mechelen/AeF6np9F/"

# Row 161 (Indice 160)
$ws.Cells.Item(161, 1).Value = 160
$ws.Cells.Item(161, 2).Value = "belgium"
$ws.Cells.Item(161, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(161, 4).Value = "2023-2024"
$ws.Cells.Item(161, 5).Value = 45287.86458333334
$ws.Cells.Item(161, 6).Value = "Anderlecht"
$ws.Cells.Item(161, 7).Value = 2
$ws.Cells.Item(161, 8).Value = "Cercle Brugge KSV"
$ws.Cells.Item(161, 9).Value = 0
$ws.Cells.Item(161, 10).Value = 1.9
$ws.Cells.Item(161, 11).Value = "23/12/2023 21:12"
$ws.Cells.Item(161, 12).Value = 2.11
$ws.Cells.Item(161, 13).Value = "27/12/2023 20:41"
$ws.Cells.Item(161, 14).Value = 3.66
$ws.Cells.Item(161, 15).Value = "23/12/2023 21:12"
$ws.Cells.Item(161, 16).Value = 3.59
$ws.Cells.Item(161, 17).Value = "27/12/2023 20:41"
$ws.Cells.Item(161, 18).Value = 3.8
$ws.Cells.Item(161, 19).Value = "23/12/2023 21:12"
$ws.Cells.Item(161, 20).Value = 3.55
$ws.Cells.Item(161, 21).Value = "27/12/2023 20:41"
$ws.Cells.Item(161, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/anderlecht-cercle-brugge/MZYBQ2nd/"
